# Generate Report for Handback
# Refresh the handback status report with the newly-processed file GUIDs,
# content hashes, and handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$newSourceGuid = "1e2b5805-6f60-4125-a897-ef3151d8ab4c"
$oldSourceGuid2 = "ffffe73ca54b-970e-4d2a-a723-8512024a563e"
$contentHash = "056ac7fcbe6e14b6529a7349561b36bd236bafa1"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newSourceGuid.md"
$wsOverview.Range("A3").Value = "$oldSourceGuid2.md"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newSourceGuid.md"
$wsZhCn.Range("D2").Value = "$newSourceGuid.$contentHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-23 05:12:15"
$wsZhCn.Range("F2").Value = "$newSourceGuid.md"
$wsZhCn.Range("G2").Value = "$newSourceGuid.$contentHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-03-23 05:12:38"

$wsZhCn.Range("A3").Value = "$oldSourceGuid2.md"
$wsZhCn.Range("D3").Value = "$newSourceGuid.$contentHash.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-23 05:12:15"
$wsZhCn.Range("F3").Value = "$oldSourceGuid2.md"
$wsZhCn.Range("G3").Value = "$newSourceGuid.$contentHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-03-23 05:12:38"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newSourceGuid.md"
$wsDeDe.Range("D2").Value = "$newSourceGuid.$contentHash.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-23 05:12:19"
$wsDeDe.Range("F2").Value = "$newSourceGuid.md"
$wsDeDe.Range("G2").Value = "$newSourceGuid.$contentHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-03-23 05:12:45"

$wsDeDe.Range("A3").Value = "$oldSourceGuid2.md"
$wsDeDe.Range("D3").Value = "$newSourceGuid.$contentHash.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-23 05:12:19"
$wsDeDe.Range("F3").Value = "$oldSourceGuid2.md"
$wsDeDe.Range("G3").Value = "$newSourceGuid.$contentHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-03-23 05:12:45"
